$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage type ("@" = Text number format) for the Price (D)
# and Hora (G) columns, matching the source data which stores these as text
# even though they look numeric.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "244.83"
$ws.Range("G2").Value = "11"

$ws.Range("D3").Value = "21.98"
$ws.Range("G3").Value = "11"

$ws.Range("D4").Value = "5.390"
$ws.Range("G4").Value = "11"

$ws.Range("D5").Value = "0.06018"
$ws.Range("G5").Value = "11"

$ws.Range("D6").Value = "3.387"
$ws.Range("G6").Value = "11"

$ws.Range("D7").Value = "0.8113"
$ws.Range("G7").Value = "11"

$ws.Range("D8").Value = "0.9564"
$ws.Range("G8").Value = "11"

$ws.Range("D9").Value = "0.1423"
$ws.Range("G9").Value = "11"

$ws.Range("D10").Value = "0.07388"
$ws.Range("G10").Value = "11"

$ws.Range("D11").Value = "0.03382"
$ws.Range("G11").Value = "11"

$ws.Range("D12").Value = "0.03052"
$ws.Range("G12").Value = "11"

$ws.Range("D13").Value = "0.09425"
$ws.Range("G13").Value = "11"

$ws.Range("D14").Value = "4.002"
$ws.Range("G14").Value = "11"

$ws.Range("G15").Value = "11"

$ws.Range("D16").Value = "0.04816"
$ws.Range("G16").Value = "11"

$ws.Range("D17").Value = "0.0005872"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("G17").Value = "11"

$ws.Range("D18").Value = "0.006268"
$ws.Range("G18").Value = "11"

$ws.Range("D19").Value = "0.005045"
$ws.Range("G19").Value = "11"

$ws.Range("D20").Value = "0.0009894"
$ws.Range("G20").Value = "11"

$ws.Range("G21").Value = "11"

$ws.Range("D22").Value = "3.697"
$ws.Range("G22").Value = "11"

$ws.Range("D23").Value = "6.426"
$ws.Range("G23").Value = "11"

$ws.Range("D24").Value = "2.185"
$ws.Range("G24").Value = "11"

$ws.Range("D25").Value = "0.3255"
$ws.Range("G25").Value = "11"

$ws.Range("D26").Value = "0.1341"
$ws.Range("G26").Value = "11"

$ws.Range("G27").Value = "11"

$ws.Range("G28").Value = "11"

$ws.Range("G29").Value = "11"

$ws.Range("G30").Value = "11"

$ws.Range("G31").Value = "11"

$ws.Range("G32").Value = "11"

$ws.Range("G33").Value = "11"

$ws.Range("G34").Value = "11"

$ws.Range("G35").Value = "11"

$ws.Range("G36").Value = "11"

$ws.Range("G37").Value = "11"

$ws.Range("G38").Value = "11"

$ws.Range("G39").Value = "11"

$ws.Range("D40").Value = "0.03991"
$ws.Range("G40").Value = "11"

$ws.Range("D41").Value = "0.006563"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "11"

$ws.Range("D42").Value = "0.1073"
$ws.Range("G42").Value = "11"

$ws.Range("D43").Value = "0.002301"
$ws.Range("G43").Value = "11"

$ws.Range("D44").Value = "0.005238"
$ws.Range("G44").Value = "11"

$ws.Range("D45").Value = "0.00005223"
$ws.Range("G45").Value = "11"

$ws.Range("G46").Value = "11"

$ws.Range("D47").Value = "0.8103"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("G47").Value = "11"

$ws.Range("D48").Value = "0.02076"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("G48").Value = "11"

$ws.Range("G49").Value = "11"

$ws.Range("G50").Value = "11"

$ws.Range("G51").Value = "11"
